$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the new columns so new cells pick up the right style
# (style "2" = font size 12 + centered, matches cols B:G already on the sheet;
#  style "1" = font size 12 only, matches col A / col K).
$ws.Range("F1:J68").Font.Size = 12
$ws.Range("F1:J68").HorizontalAlignment = -4108
$ws.Range("K1:K68").Font.Size = 12

# --- Step 1: headers for OsirisMin/OsirisMax (I/J) -- new shared strings 70,71
$ws.Range("I1").Value = "OsirisMin"
$ws.Range("J1").Value = "OsirisMax"

# --- Step 2: I/J numeric data for all rows (no new shared strings)
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 22
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 18
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 2
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 15.1
$ws.Range("I6").Value = 3.2
$ws.Range("J6").Value = 17.2
$ws.Range("I7").Value = 3.3
$ws.Range("J7").Value = 19
$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 18.2
$ws.Range("I9").Value = 7
$ws.Range("J9").Value = 26
$ws.Range("I10").Value = 10.2
$ws.Range("J10").Value = 73
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 20
$ws.Range("I12").Value = 22.3
$ws.Range("J12").Value = 43
$ws.Range("I13").Value = 2
$ws.Range("J13").Value = 29
$ws.Range("I14").Value = 1.2
$ws.Range("J14").Value = 19
$ws.Range("I15").Value = 4
$ws.Range("J15").Value = 17.1
$ws.Range("I16").Value = 11
$ws.Range("J16").Value = 28
$ws.Range("I17").Value = 5.2
$ws.Range("J17").Value = 20
$ws.Range("I24").Value = 2
$ws.Range("J24").Value = 42.3
$ws.Range("I25").Value = 0.2
$ws.Range("J25").Value = 39.2

# --- Step 3: K column -- "PP16" data first (72), then "Kit" header (73);
# the remaining PP16 cells reuse string 72, "ID" cells are deferred to step 6
# so that "ID" becomes the LAST new shared string (76), matching the source file.
$ws.Range("K2").Value = "PP16"
$ws.Range("K1").Value = "Kit"
$ws.Range("K3").Value = "PP16"
$ws.Range("K4").Value = "PP16"
$ws.Range("K5").Value = "PP16"
$ws.Range("K6").Value = "PP16"
$ws.Range("K7").Value = "PP16"
$ws.Range("K8").Value = "PP16"
$ws.Range("K9").Value = "PP16"
$ws.Range("K11").Value = "PP16"
$ws.Range("K12").Value = "PP16"
$ws.Range("K13").Value = "PP16"
$ws.Range("K14").Value = "PP16"
$ws.Range("K15").Value = "PP16"
$ws.Range("K24").Value = "PP16"
$ws.Range("K25").Value = "PP16"

# --- Step 4: headers for LadderMin/LadderMax (F/G) -- new shared strings 74,75
$ws.Range("F1").Value = "LadderMin"
$ws.Range("G1").Value = "LadderMax"

# --- Step 5: F/G numeric data for all rows
$ws.Range("F2").Value = 12
$ws.Range("G2").Value = 20
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 15
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 13.3
$ws.Range("F6").Value = 6
$ws.Range("G6").Value = 13
$ws.Range("F7").Value = 6
$ws.Range("G7").Value = 15
$ws.Range("F8").Value = 6
$ws.Range("G8").Value = 14
$ws.Range("F9").Value = 10
$ws.Range("G9").Value = 22
$ws.Range("F10").Value = 16
$ws.Range("G10").Value = 46.2
$ws.Range("F11").Value = 7
$ws.Range("G11").Value = 18
$ws.Range("F12").Value = 24
$ws.Range("G12").Value = 38
$ws.Range("F13").Value = 8
$ws.Range("G13").Value = 27
$ws.Range("F14").Value = 7
$ws.Range("G14").Value = 16
$ws.Range("F15").Value = 7
$ws.Range("G15").Value = 15
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 28
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 17.2
$ws.Range("F24").Value = 5
$ws.Range("G24").Value = 24
$ws.Range("F25").Value = 2.2
$ws.Range("G25").Value = 17

# --- Step 6: remaining K cells ("ID") -- new shared string 76 (last)
$ws.Range("K16").Value = "ID"
$ws.Range("K17").Value = "ID"

# --- Step 7: restore the selection recorded in the saved workbook
$ws.Range("A14").Select()
